$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# "Ratio" row under RETENTION -> 0.5
$c = $t.Cell(24, 4)
$c.Range.Text = "0.5"
$c.Range.Font.Bold = -1
$c.Range.Font.Size = 12
$c.Range.Font.SizeBi = 12

# "Answer Recall Lenient (ARL)" row -> 0.1666
$c = $t.Cell(44, 4)
$c.Range.Text = "0.1666"
$c.Range.Font.Bold = -1
$c.Range.Font.Size = 12
$c.Range.Font.SizeBi = 12

# "Answer Recall Strict (ARS)" row -> 0
$c = $t.Cell(45, 4)
$c.Range.Text = "0"
$c.Range.Font.Bold = -1
$c.Range.Font.Size = 12
$c.Range.Font.SizeBi = 12

# "Answer Recall Average (ARA)" row -> 0.0833
$c = $t.Cell(46, 4)
$c.Range.Text = "0.0833"
$c.Range.Font.Bold = -1
$c.Range.Font.Size = 12
$c.Range.Font.SizeBi = 12
